$p = $ppt.ActivePresentation

# --- Slide 3: "Attribute Routing" / "Explicitly declared routes"
#     becomes "Media Formatters" / "What's your vector, Victor?"
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Media Formatters"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "What’s your vector, Victor?"

# --- Slide 12: "Media Formatters" / "Built-in and custom ways to generate content"
#     becomes "Attribute Routing" / "It’s the new style to get what you want"
$s12 = $p.Slides.Item(12)
$s12.Shapes.Item(1).TextFrame.TextRange.Text = "Attribute Routing"
$s12.Shapes.Item(2).TextFrame.TextRange.Text = "It’s the new style to get what you want"

# --- Add speaker notes to slide 12 (previously had none)
$notes12 = $s12.NotesPage
$notesBody = $notes12.Shapes.AddPlaceholder(2)
$notesBody.TextFrame.TextRange.Text = "03aCustomMediaFormatter"
